# Apply the 2022FSAdates.xlsx "Add files via upload" edit.
#
# Net data changes (per the canonical-OOXML diff):
#   - Row 85 (2022-10-08, Category O, Weapon E): SplitGender 1 -> 0
#   - Row 86 (2022-10-08, Category O, Weapon S): Cancelled 0 -> 1, SplitGender 1 -> 0
#   - Row 87 (2022-10-08, Category V, Weapon F): Cancelled 0 -> 1
#   - Row 89 (2022-10-08, Category V, Weapon E): SplitGender 1 -> 0
#   - Row 90 (2022-10-08, Category V, Weapon S): Cancelled 0 -> 1
#   - Row 86, column L had a stale cached string ("20221106U13GF") left over
#     from an old formula edit; it gets restored to the live Table1 formula.
#   - The now-unused shared string "20221106U13GF" disappears from the
#     workbook once nothing references it any more (handled automatically
#     by the writer when the workbook is saved).
#   - The active sheet's view scrolls/selects a different cell
#     (was topLeftCell A107 / selection D52).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Table1 data edits (columns: A Date, B Category, C Weapon, D Cancelled,
#     E SplitGender, F Time, G Rollcall, H Name, I Link) ---

# Row 85: O / E -> no longer split by gender
$ws.Range("E85").Value = 0

# Row 86: O / S -> cancelled, and no longer split by gender
$ws.Range("D86").Value = 1
$ws.Range("E86").Value = 0

# Row 87: V / F -> cancelled
$ws.Range("D87").Value = 1

# Row 89: V / E -> no longer split by gender
$ws.Range("E89").Value = 0

# Row 90: V / S -> cancelled
$ws.Range("D90").Value = 1

# Row 86's "Womens" helper cell (L86) had been overwritten with a stale
# literal string in the source file; put the live Table1 formula back so it
# recalculates normally with the rest of the column.
$ws.Range("L86").Formula = '=IF(Table1[[#This Row],[Cancelled]]=1,"",IF(Table1[[#This Row],[SplitGender]]=0,"N/A",Table1[[#This Row],[Date]]&Table1[[#This Row],[Category]]&IF(Table1[[#This Row],[SplitGender]]=1,IF(OR(Table1[[#This Row],[Category]]="U9",Table1[[#This Row],[Category]]="U11",Table1[[#This Row],[Category]]="U13"),"G","W"),"")&Table1[[#This Row],[Weapon]]))'

# --- View state: scroll the window and move the selection ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 77
$win.ScrollColumn = 1
[void]$ws.Range("K89").Select()
